$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31 (2024-06-05)
$ws.Range("C31").Value = "fixed button add & edit component in dashboard card icon title and add menu design and conditionaly render the proper value"
$ws.Range("D31").Value = "no"
$ws.Range("E31").Value = "1 day"

# Row 32 (2024-06-06)
$ws.Range("C32").Value = "make another 4 card different design and first card in one row and other are 1 row 2 card"
$ws.Range("D32").Value = "no"
$ws.Range("E32").Value = "1 day"

# Update the active view/selection to match the author's final cursor position
$ws.Activate()
$ws.Range("C33").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
